# "updates for test problems"
# Update Inlet Mole Frac values on the "species" sheet and move that
# sheet's stored selection to the new working area (C11:C15), then
# restore focus to the originally-active "system" tab.

$wb = $excel.ActiveWorkbook
$species = $wb.Worksheets.Item("species")

$species.Range("C2").Value = 0.1
$species.Range("C4").Value = 0.4
$species.Range("C5").Value = 0.2
$species.Range("C7").Value = 0.05

[void]$species.Range("C11:C15").Select()

$wb.Worksheets.Item("system").Activate()
